$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.954100608825684
$ws.Range("B1").Value = 3.184989929199219
$ws.Range("C1").Value = 2.812205791473389
$ws.Range("D1").Value = 2.516140222549438
$ws.Range("E1").Value = 1.731523394584656
